$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.243801652892562
$ws.Range("C2").Value = 0.4586776859504132
$ws.Range("J2").Value = 0.02479338842975207
$ws.Range("P2").Value = 0.1322314049586777
$ws.Range("S2").Value = 0.140495867768595
$ws.Range("C3").Value = 0.04310344827586207
$ws.Range("J3").Value = 0.09482758620689655
$ws.Range("P3").Value = 0.5689655172413793
$ws.Range("S3").Value = 0.293103448275862
$ws.Range("J4").Value = 0.1212121212121212
$ws.Range("P4").Value = 0.4848484848484849
$ws.Range("S4").Value = 0.3939393939393939
$ws.Range("P5").Value = 1.0
$ws.Range("B6").Value = 0.0411522633744856
$ws.Range("D6").Value = 0.01646090534979424
$ws.Range("F6").Value = 0.06995884773662552
$ws.Range("J6").Value = 0.2263374485596708
$ws.Range("O6").Value = 0.0205761316872428
$ws.Range("Q6").Value = 0.1275720164609054
$ws.Range("R6").Value = 0.07407407407407407
$ws.Range("S6").Value = 0.4238683127572017
$ws.Range("B7").Value = 0.119047619047619
$ws.Range("D7").Value = 0.005952380952380952
$ws.Range("E7").Value = 0.005952380952380952
$ws.Range("F7").Value = 0.07738095238095238
$ws.Range("J7").Value = 0.125
$ws.Range("O7").Value = 0.01785714285714286
$ws.Range("Q7").Value = 0.09523809523809523
$ws.Range("R7").Value = 0.1071428571428571
$ws.Range("S7").Value = 0.4464285714285715
$ws.Range("B8").Value = 0.05555555555555555
$ws.Range("D8").Value = 0.01690821256038647
$ws.Range("F8").Value = 0.09178743961352658
$ws.Range("J8").Value = 0.1159420289855072
$ws.Range("O8").Value = 0.02898550724637681
$ws.Range("Q8").Value = 0.1473429951690821
$ws.Range("R8").Value = 0.0966183574879227
$ws.Range("S8").Value = 0.4468599033816425
$ws.Range("B9").Value = 0.1005917159763314
$ws.Range("D9").Value = 0.01775147928994083
$ws.Range("F9").Value = 0.08875739644970414
$ws.Range("J9").Value = 0.1597633136094675
$ws.Range("Q9").Value = 0.1479289940828402
$ws.Range("R9").Value = 0.1005917159763314
$ws.Range("S9").Value = 0.3846153846153846
$ws.Range("B10").Value = 0.09193548387096774
$ws.Range("D10").Value = 0.01612903225806452
$ws.Range("E10").Value = 0.0008064516129032258
$ws.Range("F10").Value = 0.0814516129032258
$ws.Range("J10").Value = 0.1040322580645161
$ws.Range("O10").Value = 0.01693548387096774
$ws.Range("Q10").Value = 0.182258064516129
$ws.Range("R10").Value = 0.1
$ws.Range("S10").Value = 0.4064516129032258
$ws.Range("G11").Value = 0.1355311355311355
$ws.Range("J11").Value = 0.09157509157509157
$ws.Range("K11").Value = 0.2051282051282051
$ws.Range("L11").Value = 0.5604395604395604
$ws.Range("S11").Value = 0.007326007326007326
$ws.Range("G12").Value = 0.75
$ws.Range("J12").Value = 0.2171052631578947
$ws.Range("S12").Value = 0.03289473684210526
$ws.Range("G13").Value = 0.6944444444444444
$ws.Range("J13").Value = 0.25
$ws.Range("S13").Value = 0.05555555555555555
$ws.Range("F15").Value = 0.02083333333333333
$ws.Range("H15").Value = 0.1833333333333333
$ws.Range("K15").Value = 0.0875
$ws.Range("M15").Value = 0.008333333333333333
$ws.Range("O15").Value = 0.04166666666666666
$ws.Range("S15").Value = 0.2083333333333333
$ws.Range("F16").Value = 0.01785714285714286
$ws.Range("H16").Value = 0.1875
$ws.Range("I16").Value = 0.07142857142857142
$ws.Range("J16").Value = 0.4464285714285715
$ws.Range("K16").Value = 0.09821428571428571
$ws.Range("M16").Value = 0.01785714285714286
$ws.Range("O16").Value = 0.03571428571428571
$ws.Range("S16").Value = 0.125
$ws.Range("F17").Value = 0.01671309192200557
$ws.Range("H17").Value = 0.1643454038997214
$ws.Range("I17").Value = 0.07799442896935933
$ws.Range("J17").Value = 0.4512534818941504
$ws.Range("K17").Value = 0.0807799442896936
$ws.Range("M17").Value = 0.02228412256267409
$ws.Range("O17").Value = 0.06685236768802229
$ws.Range("S17").Value = 0.1197771587743733
$ws.Range("F18").Value = 0.01388888888888889
$ws.Range("H18").Value = 0.1851851851851852
$ws.Range("I18").Value = 0.06018518518518518
$ws.Range("J18").Value = 0.4259259259259259
$ws.Range("K18").Value = 0.09722222222222222
$ws.Range("M18").Value = 0.02314814814814815
$ws.Range("O18").Value = 0.09722222222222222
$ws.Range("S18").Value = 0.09722222222222222
$ws.Range("F19").Value = 0.02075019952114924
$ws.Range("H19").Value = 0.1987230646448523
$ws.Range("I19").Value = 0.07980845969672785
$ws.Range("J19").Value = 0.3950518754988029
$ws.Range("K19").Value = 0.1053471667996808
$ws.Range("M19").Value = 0.01516360734237829
$ws.Range("N19").Value = 0.001596169193934557
$ws.Range("O19").Value = 0.0845969672853153
$ws.Range("S19").Value = 0.09896249002394254
